$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (Fleet Code / Agent Name / Goods & Service Tax No / Tax Type /
# Receipt Id / Receipt Date / Engine No. / Chassis No / Manufacture Year /
# Type Of Body / Type Of Fuel / Seating Capacity / Cubic Capacity added;
# "Roadtax Number" removed) - matches the export/format commit.
$headers = @(
  "Fleet Code",
  "Vehicle Number",
  "Agent Name",
  "Roadtax Amount ",
  "Goods & Service Tax No",
  "Tax Type",
  "Receipt Id",
  "Receipt Date",
  "Payment Mode",
  "Pay Number",
  "Pay Date",
  "Pay Bank",
  "Pay Branch",
  "Valid From",
  "Valid Till",
  "Engine No.",
  "Chassis No",
  "Manufacture Year",
  "Type Of Body",
  "Type Of Fuel",
  "Seating Capacity(including Driver)",
  "Cubic Capacity"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Column widths (target `width` XML attribute = ColumnWidth + 5/6; values below
# are back-solved so the exported width lands as close as possible to the
# widths recorded in the authored workbook).
$colWidths = @(
  14.165765224358974,
  16.022796474358977,
  19.59401041666667,
  21.45104166666667,
  32.73607772435898,
  12.59443108974359,
  16.879887820512824,
  18.165524839743593,
  20.1654046474359,
  19.879707532051285,
  18.736919070512823,
  16.022796474358977,
  17.73697916666667,
  20.1654046474359,
  19.59401041666667,
  23.022375801282056,
  26.879286858974364,
  25.593649839743595,
  17.308433493589746,
  20.022556089743592,
  32.73607772435898,
  24.450861378205133
)

for ($i = 0; $i -lt $colWidths.Length; $i++) {
  $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# Header row height.
$ws.Rows.Item(1).RowHeight = 15

# Selection: column D selected top to bottom, scrolled so column C is leftmost.
$ws.Range("D1:D1048576").Select()

# Workbook window geometry recorded by the author's Excel session.
$excel.Width = 13660
$excel.Height = 5120
